# M1 compatibility resolved & updated the doc
# Updates the Solana-vs-EVM comparison sheet: a handful of "Solana" column
# values were wrong/placeholder (EOSIO-era leftovers) and are corrected to
# their proper Solana terms. The cells that get real Solana content also
# lose the "TODO / needs research" red highlighting (matches the format
# already used on the neighbouring B6 cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# Row 5: "chain data storage folder" -> Solana column
#   "nodeos"  =>  "test-ledger/"
$ws.Range("B6").Copy()
$ws.Range("B5").PasteSpecial($xlPasteFormats)
$ws.Range("B5").Value = "test-ledger/"

# Row 4: "Types of accounts" -> Solana column
#   "1. User account / 2. Contract account"  =>  "1. Account / 2. Program"
$ws.Range("B6").Copy()
$ws.Range("B4").PasteSpecial($xlPasteFormats)
$ws.Range("B4").Value = "1. Account`n2. Program"

# Row 11: "32 bit" -> Solana column
#   "uint32_t"  =>  "u32"
$ws.Range("B6").Copy()
$ws.Range("B11").PasteSpecial($xlPasteFormats)
$ws.Range("B11").Value = "u32"

# Row 17: "token standard" -> Solana column
#   "EOSIO token"  =>  "Token program"  (keep the existing red "TODO" style)
$ws.Range("B17").Value = "Token program"

$excel.CutCopyMode = $false

# Update the saved cursor/selection to where the author left off editing.
$ws.Range("B18").Select()
